$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# 1) Shrink the table from 4 columns (user, password, permission, name)
#    down to 2 (user, password) by removing the "name" and "permission"
#    table columns. This keeps the ListObject/table metadata consistent
#    as we work.
$lo.ListColumns.Item(4).Delete()   # "name"
$lo.ListColumns.Item(3).Delete()   # "permission"

# 2) Physically remove the now-empty worksheet columns C:D in one shot so
#    leftover content (e.g. the styled helper cell further right) shifts
#    left by two columns, and the custom column-width definitions for the
#    old C/D columns disappear.
$ws.Range("C1:D1").EntireColumn.Delete()

# 3) Re-introduce a single new table column, "admin", holding a
#    TRUE/FALSE flag per user (replacing the old permission/name/
#    Administrador-Juan Perez style data with a boolean-ish column).
$lo.ListColumns.Add() | Out-Null
$ws.Range("C1").Value = "admin"

# Write the flag values as literal text "TRUE"/"FALSE" (not native
# booleans) by computing them with a helper formula and pasting back
# just the resulting values - this keeps the cells as shared-string
# text cells instead of boolean-typed cells.
$ws.Range("F1").Formula = '="TRUE"'
$ws.Range("F2").Formula = '="FALSE"'
$ws.Range("F3").Formula = '="FALSE"'
$ws.Range("F1:F3").Copy()
$ws.Range("C2:C4").PasteSpecial(-4163)
$excel.CutCopyMode = 0
$ws.Range("F1:F3").Clear()

# 4) Restore the usual "nothing selected / clipboard empty" UI state and
#    move the active selection to A5, matching the saved workbook view.
$ws.Range("A5").Select()
